$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object "object[,]" 18,20

$arr[0,0] = 'ECs'
$arr[0,1] = 'Efnb1'
$arr[0,2] = 'Ephb1'
$arr[0,3] = 'ECs'
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 7.102474
$arr[0,7] = 21.307422
$arr[0,8] = 0.3851819652723766
$arr[0,9] = 0.3851819652723767
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 1.397474666666667
$arr[0,13] = 4.192424
$arr[0,14] = 0.9253320219903044
$arr[0,15] = 0.9253320219903044
$arr[0,16] = 9.925527485658666
$arr[0,17] = 89.32974737092799
$arr[0,18] = 0.3564212067596875
$arr[0,19] = 0.3564212067596875

$arr[1,0] = 'ECs'
$arr[1,1] = 'Efnb1'
$arr[1,2] = 'Ephb1'
$arr[1,3] = 'MuSCs'
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 7.102474
$arr[1,7] = 21.307422
$arr[1,8] = 0.3851819652723766
$arr[1,9] = 0.3851819652723767
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 0.07803566666666666
$arr[1,13] = 0.234107
$arr[1,14] = 0.05167099121464913
$arr[1,15] = 0.05167099121464913
$arr[1,16] = 0.5542462935726666
$arr[1,17] = 4.988216642154
$arr[1,18] = 0.01990273394363026
$arr[1,19] = 0.01990273394363026

$arr[2,0] = 'ECs'
$arr[2,1] = 'Efnb1'
$arr[2,2] = 'Ephb1'
$arr[2,3] = 'Neutrophils'
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 7.102474
$arr[2,7] = 21.307422
$arr[2,8] = 0.3851819652723766
$arr[2,9] = 0.3851819652723767
$arr[2,10] = 1
$arr[2,11] = 0.3333333333333333
$arr[2,12] = 0.034731
$arr[2,13] = 0.104193
$arr[2,14] = 0.02299698679504644
$arr[2,15] = 0.02299698679504644
$arr[2,16] = 0.246676024494
$arr[2,17] = 2.220084220446
$arr[2,18] = 0.008858024569058882
$arr[2,19] = 0.008858024569058882

$arr[3,0] = 'FAPs'
$arr[3,1] = 'Efnb1'
$arr[3,2] = 'Ephb1'
$arr[3,3] = 'ECs'
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 8.299810000000001
$arr[3,7] = 24.89943
$arr[3,8] = 0.4501159915808667
$arr[3,9] = 0.4501159915808668
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 1.397474666666667
$arr[3,13] = 4.192424
$arr[3,14] = 0.9253320219903044
$arr[3,15] = 0.9253320219903044
$arr[3,16] = 11.59877421314667
$arr[3,17] = 104.38896791832
$arr[3,18] = 0.4165067406196942
$arr[3,19] = 0.4165067406196943

$arr[4,0] = 'FAPs'
$arr[4,1] = 'Efnb1'
$arr[4,2] = 'Ephb1'
$arr[4,3] = 'MuSCs'
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 8.299810000000001
$arr[4,7] = 24.89943
$arr[4,8] = 0.4501159915808667
$arr[4,9] = 0.4501159915808668
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 0.07803566666666666
$arr[4,13] = 0.234107
$arr[4,14] = 0.05167099121464913
$arr[4,15] = 0.05167099121464913
$arr[4,16] = 0.6476812065566666
$arr[4,17] = 5.82913085901
$arr[4,18] = 0.02325793944654804
$arr[4,19] = 0.02325793944654805

$arr[5,0] = 'FAPs'
$arr[5,1] = 'Efnb1'
$arr[5,2] = 'Ephb1'
$arr[5,3] = 'Neutrophils'
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 8.299810000000001
$arr[5,7] = 24.89943
$arr[5,8] = 0.4501159915808667
$arr[5,9] = 0.4501159915808668
$arr[5,10] = 1
$arr[5,11] = 0.3333333333333333
$arr[5,12] = 0.034731
$arr[5,13] = 0.104193
$arr[5,14] = 0.02299698679504644
$arr[5,15] = 0.02299698679504644
$arr[5,16] = 0.28826070111
$arr[5,17] = 2.59434630999
$arr[5,18] = 0.01035131151462443
$arr[5,19] = 0.01035131151462443

$arr[6,0] = 'Inflammatory-Mac'
$arr[6,1] = 'Efnb1'
$arr[6,2] = 'Ephb1'
$arr[6,3] = 'ECs'
$arr[6,4] = 2
$arr[6,5] = 0.6666666666666666
$arr[6,6] = 0.2555593333333333
$arr[6,7] = 0.766678
$arr[6,8] = 0.01385951518541732
$arr[6,9] = 0.01385951518541733
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 1.397474666666667
$arr[6,13] = 4.192424
$arr[6,14] = 0.9253320219903044
$arr[6,15] = 0.9253320219903044
$arr[6,16] = 0.3571376941635555
$arr[6,17] = 3.214239247472
$arr[6,18] = 0.01282465321032754
$arr[6,19] = 0.01282465321032754

$arr[7,0] = 'Inflammatory-Mac'
$arr[7,1] = 'Efnb1'
$arr[7,2] = 'Ephb1'
$arr[7,3] = 'MuSCs'
$arr[7,4] = 2
$arr[7,5] = 0.6666666666666666
$arr[7,6] = 0.2555593333333333
$arr[7,7] = 0.766678
$arr[7,8] = 0.01385951518541732
$arr[7,9] = 0.01385951518541733
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 0.07803566666666666
$arr[7,13] = 0.234107
$arr[7,14] = 0.05167099121464913
$arr[7,15] = 0.05167099121464913
$arr[7,16] = 0.01994274294955555
$arr[7,17] = 0.179484686546
$arr[7,18] = 0.0007161348873849947
$arr[7,19] = 0.0007161348873849949

$arr[8,0] = 'Inflammatory-Mac'
$arr[8,1] = 'Efnb1'
$arr[8,2] = 'Ephb1'
$arr[8,3] = 'Neutrophils'
$arr[8,4] = 2
$arr[8,5] = 0.6666666666666666
$arr[8,6] = 0.2555593333333333
$arr[8,7] = 0.766678
$arr[8,8] = 0.01385951518541732
$arr[8,9] = 0.01385951518541733
$arr[8,10] = 1
$arr[8,11] = 0.3333333333333333
$arr[8,12] = 0.034731
$arr[8,13] = 0.104193
$arr[8,14] = 0.02299698679504644
$arr[8,15] = 0.02299698679504644
$arr[8,16] = 0.008875831205999999
$arr[8,17] = 0.07988248085399999
$arr[8,18] = 0.0003187270877047878
$arr[8,19] = 0.0003187270877047879

$arr[9,0] = 'MuSCs'
$arr[9,1] = 'Efnb1'
$arr[9,2] = 'Ephb1'
$arr[9,3] = 'ECs'
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 1.392645666666667
$arr[9,7] = 4.177937
$arr[9,8] = 0.07552607652132563
$arr[9,9] = 0.07552607652132566
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 1.397474666666667
$arr[9,13] = 4.192424
$arr[9,14] = 0.9253320219903044
$arr[9,15] = 0.9253320219903044
$arr[9,16] = 1.946187038809778
$arr[9,17] = 17.515683349288
$arr[9,18] = 0.0698866971004727
$arr[9,19] = 0.06988669710047272

$arr[10,0] = 'MuSCs'
$arr[10,1] = 'Efnb1'
$arr[10,2] = 'Ephb1'
$arr[10,3] = 'MuSCs'
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 1.392645666666667
$arr[10,7] = 4.177937
$arr[10,8] = 0.07552607652132563
$arr[10,9] = 0.07552607652132566
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 0.07803566666666666
$arr[10,13] = 0.234107
$arr[10,14] = 0.05167099121464913
$arr[10,15] = 0.05167099121464913
$arr[10,16] = 0.1086760330287778
$arr[10,17] = 0.9780842972589999
$arr[10,18] = 0.003902507236410335
$arr[10,19] = 0.003902507236410336

$arr[11,0] = 'MuSCs'
$arr[11,1] = 'Efnb1'
$arr[11,2] = 'Ephb1'
$arr[11,3] = 'Neutrophils'
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 1.392645666666667
$arr[11,7] = 4.177937
$arr[11,8] = 0.07552607652132563
$arr[11,9] = 0.07552607652132566
$arr[11,10] = 1
$arr[11,11] = 0.3333333333333333
$arr[11,12] = 0.034731
$arr[11,13] = 0.104193
$arr[11,14] = 0.02299698679504644
$arr[11,15] = 0.02299698679504644
$arr[11,16] = 0.048367976649
$arr[11,17] = 0.435311789841
$arr[11,18] = 0.001736872184442593
$arr[11,19] = 0.001736872184442593

$arr[12,0] = 'Neutrophils'
$arr[12,1] = 'Efnb1'
$arr[12,2] = 'Ephb1'
$arr[12,3] = 'ECs'
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 1.139971333333333
$arr[12,7] = 3.419914
$arr[12,8] = 0.06182302089771886
$arr[12,9] = 0.06182302089771888
$arr[12,10] = 3
$arr[12,11] = 1
$arr[12,12] = 1.397474666666667
$arr[12,13] = 4.192424
$arr[12,14] = 0.9253320219903044
$arr[12,15] = 0.9253320219903044
$arr[12,16] = 1.593081059059556
$arr[12,17] = 14.337729531536
$arr[12,18] = 0.05720682093283504
$arr[12,19] = 0.05720682093283506

$arr[13,0] = 'Neutrophils'
$arr[13,1] = 'Efnb1'
$arr[13,2] = 'Ephb1'
$arr[13,3] = 'MuSCs'
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 1.139971333333333
$arr[13,7] = 3.419914
$arr[13,8] = 0.06182302089771886
$arr[13,9] = 0.06182302089771888
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 0.07803566666666666
$arr[13,13] = 0.234107
$arr[13,14] = 0.05167099121464913
$arr[13,15] = 0.05167099121464913
$arr[13,16] = 0.08895842297755555
$arr[13,17] = 0.800625806798
$arr[13,18] = 0.003194456769669101
$arr[13,19] = 0.003194456769669102

$arr[14,0] = 'Neutrophils'
$arr[14,1] = 'Efnb1'
$arr[14,2] = 'Ephb1'
$arr[14,3] = 'Neutrophils'
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 1.139971333333333
$arr[14,7] = 3.419914
$arr[14,8] = 0.06182302089771886
$arr[14,9] = 0.06182302089771888
$arr[14,10] = 1
$arr[14,11] = 0.3333333333333333
$arr[14,12] = 0.034731
$arr[14,13] = 0.104193
$arr[14,14] = 0.02299698679504644
$arr[14,15] = 0.02299698679504644
$arr[14,16] = 0.039592344378
$arr[14,17] = 0.356331099402
$arr[14,18] = 0.001421743195214721
$arr[14,19] = 0.001421743195214721

$arr[15,0] = 'Resolving-Mac'
$arr[15,1] = 'Efnb1'
$arr[15,2] = 'Ephb1'
$arr[15,3] = 'ECs'
$arr[15,4] = 2
$arr[15,5] = 0.6666666666666666
$arr[15,6] = 0.248809
$arr[15,7] = 0.746427
$arr[15,8] = 0.0134934305422948
$arr[15,9] = 0.01349343054229481
$arr[15,10] = 3
$arr[15,11] = 1
$arr[15,12] = 1.397474666666667
$arr[15,13] = 4.192424
$arr[15,14] = 0.9253320219903044
$arr[15,15] = 0.9253320219903044
$arr[15,16] = 0.3477042743386666
$arr[15,17] = 3.129338469048
$arr[15,18] = 0.01248590336728738
$arr[15,19] = 0.01248590336728738

$arr[16,0] = 'Resolving-Mac'
$arr[16,1] = 'Efnb1'
$arr[16,2] = 'Ephb1'
$arr[16,3] = 'MuSCs'
$arr[16,4] = 2
$arr[16,5] = 0.6666666666666666
$arr[16,6] = 0.248809
$arr[16,7] = 0.746427
$arr[16,8] = 0.0134934305422948
$arr[16,9] = 0.01349343054229481
$arr[16,10] = 3
$arr[16,11] = 1
$arr[16,12] = 0.07803566666666666
$arr[16,13] = 0.234107
$arr[16,14] = 0.05167099121464913
$arr[16,15] = 0.05167099121464913
$arr[16,16] = 0.01941597618766666
$arr[16,17] = 0.174743785689
$arr[16,18] = 0.0006972189310063931
$arr[16,19] = 0.0006972189310063933

$arr[17,0] = 'Resolving-Mac'
$arr[17,1] = 'Efnb1'
$arr[17,2] = 'Ephb1'
$arr[17,3] = 'Neutrophils'
$arr[17,4] = 2
$arr[17,5] = 0.6666666666666666
$arr[17,6] = 0.248809
$arr[17,7] = 0.746427
$arr[17,8] = 0.0134934305422948
$arr[17,9] = 0.01349343054229481
$arr[17,10] = 1
$arr[17,11] = 0.3333333333333333
$arr[17,12] = 0.034731
$arr[17,13] = 0.104193
$arr[17,14] = 0.02299698679504644
$arr[17,15] = 0.02299698679504644
$arr[17,16] = 0.008641385378999999
$arr[17,17] = 0.077772468411
$arr[17,18] = 0.00031030824400103
$arr[17,19] = 0.00031030824400103

$ws.Range("A2:T19").Value2 = $arr

Write-Output "done"